$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace Variable1/Wert1/Value1 -> BaseStateChromium/Started/Not Started
$ws.Range("A2").Value = "BaseStateChromium"
$ws.Range("B2").Value = "Started"
$ws.Range("C2").Value = "Not Started"

# Row 3: replace Variable2/Wert2/Value2 -> BaseStatePixel9Pro_API35/Not Started/Started
$ws.Range("A3").Value = "BaseStatePixel9Pro_API35"
$ws.Range("B3").Value = "Not Started"
$ws.Range("C3").Value = "Started"

# Row 4: new row AUT / Chromium / Pixel9Pro_API35
$ws.Range("A4").Value = "AUT"
$ws.Range("B4").Value = "Chromium"
$ws.Range("C4").Value = "Pixel9Pro_API35"

# Update selection to match target view state
$ws.Range("F20").Select()
